# Sheet1 holds the TestNG test-case table. The old first column (TC_ID /
# TC_01) is being repurposed into a "Functions" column that names the
# actual test method (TC001_TestNG_POC) to run via the Google-search test
# script referenced in the commit message.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = "Functions"
$ws.Range("A2").Value = "TC001_TestNG_POC"

# Widen column A so the new, longer header/value are fully visible
# (~23.29 stored width == 22.57 characters at the default Calibri 11 MDW).
$ws.Columns.Item(1).ColumnWidth = 22.57

# Move the active selection off the hyperlink cell and onto the new column.
$ws.Range("A2").Select() | Out-Null
